$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.964.36"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "3.239.01"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'580.24"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").Value = "'184.37"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "'0.129"
$ws.Range("E9").Value = "  -3.75%  "

$ws.Range("E10").Value = "  -1.11%  "

$ws.Range("D11").Value = "'0.414"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "3.804.61"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("D14").Value = "'27.71"
$ws.Range("E14").Value = "  -3.30%  "

$ws.Range("D15").Value = "67.981.70"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").Value = "3.237.58"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "'5.78"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").Value = "'396.22"
$ws.Range("E20").Value = "  +4.12%  "

$ws.Range("D21").Value = "'7.56"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'71.16"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'0.514"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("E25").Value = "  -1.27%  "

$ws.Range("D26").Value = "'0.186"
$ws.Range("E26").Value = "  +2.75%  "

$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").Value = "'5.58"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("D31").Value = "'22.72"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").Value = "'7.01"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").Value = "'1.25"
$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "'162.08"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "  -4.44%  "

$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("D38").Value = "'26.51"
$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").Value = "'0.810"
$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").Value = "'4.56"
$ws.Range("E40").Value = "  -0.57%  "

$ws.Range("D41").Value = "'6.48"
$ws.Range("E41").Value = "  -3.28%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'24.95"
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("D46").Value = "2.602.77"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").Value = "'336.02"
$ws.Range("E47").Value = "  -2.91%  "

$ws.Range("D48").Value = "'0.0278"
$ws.Range("E48").Value = "  -1.79%  "

$ws.Range("D49").Value = "'6.29"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("E50").Value = "  -1.58%  "

$ws.Range("D51").Value = "'31.05"
$ws.Range("E51").Value = "  +2.07%  "

# Row 42/43: OKB and dogwifhat swap positions with updated values
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'41.21"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.47"
$ws.Range("E43").Value = "  -4.45%  "
